# Cover letter update: the manuscript is now being submitted to
# "American Naturalist" instead of "Evolution" ("updated for am nat").
#
# This touches a single paragraph near the end of the letter:
#   "We hope you will find this manuscript worthy of publication in
#    Evolution."
# becomes
#   "We hope you will find this manuscript worthy of publication in
#    American Naturalist."
#
# The paragraph also carries Word's "last edit position" bookmark
# (_GoBack). In the target revision that bookmark sits right after the
# journal name (i.e. where the user's cursor was after typing the
# replacement), instead of at the very start of the paragraph. We
# recreate it in that spot after performing the text swap.

$d = $word.ActiveDocument

$oldJournal = "Evolution"
$newJournal = "American Naturalist"

# --- Step 1: drop the existing _GoBack bookmark (if any). We'll add it
#     back after the text edit, positioned right after the new journal
#     name, so it doesn't interfere with locating/replacing the text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: locate the journal name and replace it. Use the plain
#     document text (case-sensitive) so we don't accidentally match the
#     lower-case "evolution" that appears earlier in the letter.
$text = $d.Content.Text
$start = $text.IndexOf($oldJournal)
if ($start -lt 0) {
    throw "Could not find '$oldJournal' in the document"
}

$target = $d.Range($start, $start + $oldJournal.Length)
$target.Text = $newJournal

# --- Step 3: re-insert the _GoBack bookmark immediately after the new
#     journal name (collapsed range, i.e. a simple cursor position).
$text2 = $d.Content.Text
$bmPos = $text2.IndexOf($newJournal) + $newJournal.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
